# Updates cryptos list figures (price/volume columns) to match the
# latest scrape, and fixes the WrappedBTC/WrappedEther row ordering
# (rows 17-18 were swapped upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.127.51'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '3.513.40'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.69'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.48'
$ws.Range("E6").Value = '  -3.85%  '
$ws.Range("D7").Value = '3.504.93'
$ws.Range("E7").Value = '  -2.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("E8").Value = '  -3.08%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  +1.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.653'
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.30'
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000300'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.45'
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("D15").Value = '4.077.12'
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.44'
$ws.Range("E16").Value = '  -2.61%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.518.40'
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.094.36'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.32'
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.119'
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '544.09'
$ws.Range("E21").Value = '  +14.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.01'
$ws.Range("E22").Value = '  -3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.37'
$ws.Range("E23").Value = '  -6.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.00'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.42'
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.33'
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.04'
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").Value = '  -2.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.75'
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.27'
$ws.Range("E31").Value = '  -5.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.67'
$ws.Range("E32").Value = '  +3.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.64'
$ws.Range("E33").Value = '  -2.74%  '
$ws.Range("E34").Value = '  -4.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '556.19'
$ws.Range("E35").Value = '  -5.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.10'
$ws.Range("E36").Value = '  +7.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.03'
$ws.Range("E37").Value = '  -2.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.401'
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  -5.11%  '
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("E42").Value = '  -4.99%  '
$ws.Range("E43").Value = '  -3.70%  '
$ws.Range("D44").Value = '3.280.12'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("E45").Value = '  -2.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0445'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.47'
$ws.Range("E47").Value = '  +2.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.134'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.88'
$ws.Range("E49").Value = '  -6.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.59'
$ws.Range("E51").Value = '  +1.94%  '
